$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row2
$ws.Range("G2").Value = 0.7310083333333334
$ws.Range("H2").Value = 2.193025
$ws.Range("I2").Value = 0.01673731480740535
$ws.Range("J2").Value = 0.01673731480740535
$ws.Range("M2").Value = 14.48297233333333
$ws.Range("N2").Value = 43.448917
$ws.Range("O2").Value = 0.2019336017030403
$ws.Range("P2").Value = 0.2019336017030403
$ws.Range("Q2").Value = 10.58717346710278
$ws.Range("R2").Value = 95.28456120392501
$ws.Range("S2").Value = 0.003379826261896991
$ws.Range("T2").Value = 0.003379826261896991

# row3
$ws.Range("G3").Value = 0.7310083333333334
$ws.Range("H3").Value = 2.193025
$ws.Range("I3").Value = 0.01673731480740535
$ws.Range("J3").Value = 0.01673731480740535
$ws.Range("O3").Value = 0.007144147385663391
$ws.Range("P3").Value = 0.00714414738566339
$ws.Range("Q3").Value = 0.374560385238889
$ws.Range("R3").Value = 3.37104346715
$ws.Range("S3").Value = 0.0001195738438243501
$ws.Range("T3").Value = 0.0001195738438243501

# row4
$ws.Range("G4").Value = 0.7310083333333334
$ws.Range("H4").Value = 2.193025
$ws.Range("I4").Value = 0.01673731480740535
$ws.Range("J4").Value = 0.01673731480740535
$ws.Range("M4").Value = 56.726097
$ws.Range("N4").Value = 170.178291
$ws.Range("O4").Value = 0.7909222509112964
$ws.Range("P4").Value = 0.7909222509112963
$ws.Range("Q4").Value = 41.467249624475
$ws.Range("R4").Value = 373.205246620275
$ws.Range("S4").Value = 0.01323791470168401
$ws.Range("T4").Value = 0.01323791470168401

# row5
$ws.Range("I5").Value = 0.8536212576586365
$ws.Range("J5").Value = 0.8536212576586365
$ws.Range("M5").Value = 14.48297233333333
$ws.Range("N5").Value = 43.448917
$ws.Range("O5").Value = 0.2019336017030403
$ws.Range("P5").Value = 0.2019336017030403
$ws.Range("Q5").Value = 539.9573607852467
$ws.Range("R5").Value = 4859.61624706722
$ws.Range("S5").Value = 0.1723748150492874
$ws.Range("T5").Value = 0.1723748150492874

# row6
$ws.Range("I6").Value = 0.8536212576586365
$ws.Range("J6").Value = 0.8536212576586365
$ws.Range("O6").Value = 0.007144147385663391
$ws.Range("P6").Value = 0.00714414738566339
$ws.Range("Q6").Value = 19.10298699617333
$ws.Range("S6").Value = 0.006098396076248644
$ws.Range("T6").Value = 0.006098396076248643

# row7
$ws.Range("I7").Value = 0.8536212576586365
$ws.Range("J7").Value = 0.8536212576586365
$ws.Range("M7").Value = 56.726097
$ws.Range("N7").Value = 170.178291
$ws.Range("O7").Value = 0.7909222509112964
$ws.Range("P7").Value = 0.7909222509112963
$ws.Range("Q7").Value = 2114.87482809534
$ws.Range("R7").Value = 19033.87345285806
$ws.Range("S7").Value = 0.6751480465331005
$ws.Range("T7").Value = 0.6751480465331005

# row8
$ws.Range("G8").Value = 5.662136666666666
$ws.Range("H8").Value = 16.98641
$ws.Range("I8").Value = 0.129641427533958
$ws.Range("J8").Value = 0.129641427533958
$ws.Range("M8").Value = 14.48297233333333
$ws.Range("N8").Value = 43.448917
$ws.Range("O8").Value = 0.2019336017030403
$ws.Range("P8").Value = 0.2019336017030403
$ws.Range("Q8").Value = 82.00456869088555
$ws.Range("R8").Value = 738.0411182179701
$ws.Range("S8").Value = 0.02617896039185584
$ws.Range("T8").Value = 0.02617896039185584

# row9
$ws.Range("G9").Value = 5.662136666666666
$ws.Range("H9").Value = 16.98641
$ws.Range("I9").Value = 0.129641427533958
$ws.Range("J9").Value = 0.129641427533958
$ws.Range("O9").Value = 0.007144147385663391
$ws.Range("P9").Value = 0.00714414738566339
$ws.Range("Q9").Value = 2.901214657117778
$ws.Range("R9").Value = 26.11093191406
$ws.Range("S9").Value = 0.0009261774655903961
$ws.Range("T9").Value = 0.000926177465590396

# row10
$ws.Range("G10").Value = 5.662136666666666
$ws.Range("H10").Value = 16.98641
$ws.Range("I10").Value = 0.129641427533958
$ws.Range("J10").Value = 0.129641427533958
$ws.Range("M10").Value = 56.726097
$ws.Range("N10").Value = 170.178291
$ws.Range("O10").Value = 0.7909222509112964
$ws.Range("P10").Value = 0.7909222509112963
$ws.Range("Q10").Value = 321.19091378059
$ws.Range("R10").Value = 2890.71822402531
$ws.Range("S10").Value = 0.1025362896765118
$ws.Range("T10").Value = 0.1025362896765118

